$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.925.40"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.904.09"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.8052"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +6.81%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "241.19"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.84%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.9997"
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3117"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +2.76%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "26.40"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +4.68%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07010"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +3.22%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07995"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.908.91"
$ws.Range("E12").Value = "  +0.63%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.7426"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +0.30%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.179"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.81%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "92.34"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "29.935.81"
$ws.Range("E16").Value = "  +0.59%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.96"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +0.90%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "5.861"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.54%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "245.08"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.15%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007775"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.93%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "2.152.84"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  -0.11%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.925"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.90%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "168.17"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +1.75%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.203"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.25%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1475"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +16.21%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.86"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.38%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.062"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +2.86%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.361"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("E31").Value = "  -0.07%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.291"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +1.55%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05528"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +5.93%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.063"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.10%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.260"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +1.44%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7283"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +0.74%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.713"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +0.19%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01916"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.99%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.782"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +0.44%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.4399"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +0.67%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "71.98"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +1.11%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.963"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -2.47%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.9991"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.17%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8369"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +1.60%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.885"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +0.60%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "100.75"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +1.27%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "7.560"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +0.48%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.702"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "981.91"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +9.60%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.058.46"
$ws.Range("E50").Value = "  +0.40%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "36.14"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.97%  "
